$d = $word.ActiveDocument

function Insert-Fragment($range, $bodyXml) {
    $pkg = '<?xml version="1.0" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' + $bodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# --- 1) "Mohammad Aljagthmi" -> split with proofErr around the surname ---
$r1 = $d.Content
$r1.Find.Execute("Mohammad Aljagthmi") | Out-Null
$body1 = '<w:p w14:paraId="76493A8D" w14:textId="77777777" w:rsidR="009C4EA6" w:rsidRPr="008F0281" w:rsidRDefault="009C4EA6" w:rsidP="00A46CD5">' +
  '<w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">Mohammad </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Aljagthmi</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '</w:p>'
Insert-Fragment $r1 $body1

# --- 2) "Jake Manser" -> split with proofErr around the surname ---
$r2 = $d.Content
$r2.Find.Execute("Jake Manser") | Out-Null
$body2 = '<w:p w14:paraId="582C2D1E" w14:textId="77777777" w:rsidR="009C4EA6" w:rsidRPr="008F0281" w:rsidRDefault="009C4EA6" w:rsidP="00A46CD5">' +
  '<w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">Jake </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Manser</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '</w:p>'
Insert-Fragment $r2 $body2

# --- 3) Split the long "household tracker" paragraph around "WiFi", and merge the
#        following empty bookmark-only paragraph into it (deleting that blank line) ---
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "This document describes a household tracker*") {
        $target = $p
        break
    }
}
$nextPara = $target.Next()
$r3 = $d.Range($target.Range.Start, $nextPara.Range.End)
$body3 = '<w:p w14:paraId="6C965C9E" w14:textId="0A7DCB7D" w:rsidR="009C4EA6" w:rsidRDefault="009C4EA6" w:rsidP="009C4EA6">' +
  '<w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr>' +
  '<w:r w:rsidRPr="00283F70"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr>' +
  '<w:t xml:space="preserve">This document describes a household tracker that can keep track of the location of household objects. The tracker can be attached and removed from household objects manually and without extra tools. Tracker connects and syncs with application using </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>WiFi</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> to interface with the user.</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '</w:p>'
Insert-Fragment $r3 $body3

Write-Output "applied edits"
